$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "51.790.12", "  +0.33%  ")
  3 = @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.016.87", "  +2.76%  ")
  4 = @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.999", "  -0.08%  ")
  5 = @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "381.26", "  +4.98%  ")
  6 = @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "105.95", "  +0.15%  ")
  7 = @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.546", "  -0.26%  ")
  8 = @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.01%  ")
  9 = @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.605", "  +1.32%  ")
  10 = @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "37.98", "  +1.65%  ")
  11 = @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.140", "  +0.08%  ")
  12 = @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0850", "  +0.91%  ")
  13 = @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "18.91", "  +0.85%  ")
  14 = @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.491.18", "  +2.57%  ")
  15 = @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "7.57", "  +1.98%  ")
  16 = @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.017.31", "  +2.25%  ")
  17 = @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.987", "  +1.42%  ")
  18 = @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "51.843.44", "  +0.59%  ")
  19 = @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "3.49", "  +4.52%  ")
  20 = @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.50", "  +2.37%  ")
  21 = @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "13.19", "  -0.04%  ")
  22 = @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0966", "  +1.25%  ")
  23 = @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "69.10", "  +0.16%  ")
  24 = @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "265.16", "  +0.47%  ")
  25 = @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.80", "  +3.04%  ")
  26 = @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "7.44", "  +18.68%  ")
  27 = @("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.174", "  -1.49%  ")
  28 = @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "26.30", "  -0.32%  ")
  29 = @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "7.47", "  +1.78%  ")
  30 = @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.06%  ")
  31 = @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.106", "  -4.75%  ")
  32 = @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.97", "  -1.24%  ")
  33 = @("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "34.82", "  -1.23%  ")
  34 = @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "51.35", "  -0.59%  ")
  35 = @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.08", "  -3.74%  ")
  36 = @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0448", "  +4.68%  ")
  37 = @("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  +0.03%  ")
  38 = @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.15", "  -1.86%  ")
  39 = @("Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "17.60", "  +2.09%  ")
  40 = @("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "2.67", "  -6.37%  ")
  41 = @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.88", "  -0.68%  ")
  42 = @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.117", "  +1.83%  ")
  43 = @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "124.83", "  +3.77%  ")
  44 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "22.54", "  -1.41%  ")
  45 = @("WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.09", "  -2.75%  ")
  46 = @("ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "2.46", "  +7.13%  ")
  47 = @("TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.278", "  +16.64%  ")
  48 = @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.058.07", "  -1.73%  ")
  49 = @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.29", "  +1.56%  ")
  50 = @("BEAM", "https://coinranking.com/coin/cYYMfXF4u+beam-beam", "0.0349", "  +10.04%  ")
  51 = @("SEI", "https://coinranking.com/coin/8nxCqs-uj+sei-sei", "0.876", "  -0.04%  ")
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $rng = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 5))
  $rng.NumberFormat = "@"
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).Value = $vals[2]
  $ws.Cells.Item($row, 5).Value = $vals[3]
}
